$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2, $firstCol, $lastCol) {
    # Capture all values from both rows first (column A / index col is NOT included,
    # since it must stay fixed in place).
    $vals1 = @{}
    $vals2 = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals1[$c] = $ws.Cells.Item($r1, $c).Value2
        $vals2[$c] = $ws.Cells.Item($r2, $c).Value2
    }

    # Write row2's original values into row1, and row1's original values into row2.
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        if ($null -eq $vals2[$c]) {
            $cell1.ClearContents()
        } else {
            $cell1.Value2 = $vals2[$c]
        }

        if ($null -eq $vals1[$c]) {
            $cell2.ClearContents()
        } else {
            $cell2.Value2 = $vals1[$c]
        }
    }
}

# Columns B (2) through AD (30) hold the match data; column A (1) is the fixed
# row index and must not move.
Swap-RowData 52 53 2 30
Swap-RowData 130 131 2 30
Swap-RowData 167 168 2 30
